$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.617890119552612
$ws.Range("B1").Value = 3.28220796585083
$ws.Range("C1").Value = 4.333929538726807
$ws.Range("D1").Value = 1.342411398887634
$ws.Range("E1").Value = 0.7859049439430237
